# "Add files via upload" - refreshed source data (Sheet3 lookup table) plus a
# new daily snapshot column ("09-nov") appended to Sheet1, one column to the
# right of the previous last date column ("08-nov" / CK).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 1) Refresh the raw lookup values in Sheet3!B20:B36 (source data behind the
#    VLOOKUP in Sheet3!C2:C18 and in Sheet1's CB/CC formula columns).
#    Rows whose value is unchanged in the new upload (B21, B22, B25) are left
#    untouched.
# ---------------------------------------------------------------------------
$ws3.Range("B20").Value = 13.882183596107454   # 3D QUESO 92GX27
$ws3.Range("B23").Value = 6.8770735298167072   # DORITOS QUESO 70X40G
$ws3.Range("B24").Value = 6.3560404930784555   # DORITOS QUESO 77GX26
$ws3.Range("B26").Value = 4.9656565742334378   # LAYS CLASICAS 145GRX18
$ws3.Range("B27").Value = 14.149124892046276   # LAYS CLASICAS 249GRX14
$ws3.Range("B28").Value = 4.7401233884939167   # LAYS CLASICAS 40GX68
$ws3.Range("B29").Value = 4.7740715805220697   # LAYS CLASICAS 94GRX25
$ws3.Range("B30").Value = 2.9760059985002347   # LAYS ONDAS FH 30GX72
$ws3.Range("B31").Value = 11.126730929235961   # LAYS ONDAS FH 70GX28
$ws3.Range("B32").Value = 3.0404779272033577   # LAYS QSO Y CEBOLLA 34GX72
$ws3.Range("B33").Value = 17.772727760275266   # PEHUAMAR ACANALADA 520GX9
$ws3.Range("B34").Value = 9.5280916333972634   # PEHUAMAR MAICITOS 285GX10
$ws3.Range("B35").Value = 5.463719645797485    # PEHUAMAR PAPA LISA 520GX9
$ws3.Range("B36").Value = 47.775264882556257   # QUAKER AVENA INSTANT FORTIF 18X280G

# Sheet3!C2:C18 (IFERROR/VLOOKUP against A20:B36) and Sheet1!CB:CC (VLOOKUP
# against Sheet3!B1:C18) recalc automatically from the edits above.

# ---------------------------------------------------------------------------
# 2) Append the new daily snapshot column CL ("09-nov") to Sheet1, right
#    after the previous last column CK ("08-nov"). Header + 17 data rows,
#    mirroring CK's number format/style and copying the now-refreshed CB
#    value (today's VLOOKUP result) into each row as a static snapshot.
# ---------------------------------------------------------------------------
$ws1.Range("CL1").NumberFormat = $ws1.Range("CK1").NumberFormat
$ws1.Range("CL1").Value = "09-nov"

for ($r = 2; $r -le 18; $r++) {
    $cl = $ws1.Cells.Item($r, 90)   # column CL
    $cb = $ws1.Cells.Item($r, 80)   # column CB
    $cl.NumberFormat = $ws1.Cells.Item($r, 89).NumberFormat   # mirror CK's format
    $cl.Value = $cb.Value2
}

# ---------------------------------------------------------------------------
# 3) Match the author's final selection (cell CN9 on Sheet1).
# ---------------------------------------------------------------------------
$null = $ws1.Range("CN9").Select()
